# Add "Wins", "Losses", "Ties" columns (AD, AE, AF) to the worksheet.
# Header row (row 1) gets the column titles (with the same bold/centered
# header style used by the rest of row 1); every data row (2-51) gets the
# season record values: Wins=81, Losses=81, Ties=0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (from AC1) onto the three new header
# cells so they match the rest of row 1's formatting.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Find the last used row (data rows 2..51) and fill in the season record.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 81   # AD = column 30 -> Wins
    $ws.Cells.Item($r, 31).Value = 81   # AE = column 31 -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF = column 32 -> Ties
}
